# Update countries & provincias Spain
# Refresh the COVID-19 "Pais" sheet with the latest daily snapshot:
#   - update the "last refreshed" timestamp in A1
#   - update per-country stats (Casos totales/Nuevos casos/Casos activos/
#     Recuperados/Casos criticos/Muertes hoy/Muertes) for the countries whose
#     numbers moved in the new pull, including a handful of countries that
#     swapped rows because the table is sorted by 'Casos totales' descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos actualizados ... (header timestamp)
$ws.Range("A1").Value = 'Datos actualizados a 2 de Septiembre de 2020 a las 00:14'

# Row 4: Estados Unidos
$ws.Range("B4").Value = 6252488
$ws.Range("C4").Value = 36896
$ws.Range("D4").Value = 3480623
$ws.Range("E4").Value = 2583156
$ws.Range("G4").Value = 973
$ws.Range("H4").Value = 188709

# Row 5: Brasil
$ws.Range("B5").Value = 3950931
$ws.Range("C5").Value = 40030
$ws.Range("E5").Value = 730601
$ws.Range("G5").Value = 1081
$ws.Range("H5").Value = 122596

# Row 8: Mexico
$ws.Range("B8").Value = 657129
$ws.Range("C8").Value = 5092
$ws.Range("D8").Value = 471599
$ws.Range("E8").Value = 156462
$ws.Range("G8").Value = 124
$ws.Range("H8").Value = 29068

# Row 30: Canada
$ws.Range("B30").Value = 118538
$ws.Range("C30").Value = 1942
$ws.Range("D30").Value = 96214
$ws.Range("E30").Value = 21367
$ws.Range("G30").Value = 18
$ws.Range("H30").Value = 957

# Row 34: Bolivia
$ws.Range("B34").Value = 99115
$ws.Range("C34").Value = 176
$ws.Range("D34").Value = 73828
$ws.Range("E34").Value = 19847
$ws.Range("G34").Value = 19
$ws.Range("H34").Value = 5440

# Row 79: Bosnia y Herzegovina
$ws.Range("B79").Value = 18103
$ws.Range("C79").Value = 36
$ws.Range("D79").Value = 16814
$ws.Range("E79").Value = 1172

# Row 83: Libano
$ws.Range("B83").Value = 16454
$ws.Range("C83").Value = 188
$ws.Range("D83").Value = 11615
$ws.Range("E83").Value = 4197
$ws.Range("G83").Value = 13
$ws.Range("H83").Value = 642

# Row 105: Luxemburgo
$ws.Range("B105").Value = 6702
$ws.Range("C105").Value = 25

# Row 108: Malaui
$ws.Range("B108").Value = 5576
$ws.Range("C108").Value = 10
$ws.Range("D108").Value = 3420
$ws.Range("E108").Value = 1981

# Row 116: Ruanda
$ws.Range("A116").Value = 'Ruanda'
$ws.Range("B116").Value = 4142
$ws.Range("C116").Value = 79
$ws.Range("D116").Value = 2044
$ws.Range("E116").Value = 2082
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 16

# Row 117: Cuba
$ws.Range("A117").Value = 'Cuba'
$ws.Range("B117").Value = 4065
$ws.Range("C117").Value = 33
$ws.Range("D117").Value = 3395
$ws.Range("E117").Value = 575
$ws.Range("G117").Value = 1
$ws.Range("H117").Value = 95

# Row 143: Sierra Leona
$ws.Range("B143").Value = 2028
$ws.Range("C143").Value = 6
$ws.Range("E143").Value = 363
$ws.Range("G143").Value = 1
$ws.Range("H143").Value = 71

# Row 145: Yemen
$ws.Range("B145").Value = 1962
$ws.Range("C145").Value = 4
$ws.Range("D145").Value = 1160
$ws.Range("E145").Value = 232
$ws.Range("G145").Value = 4
$ws.Range("H145").Value = 570

# Row 153: Republica de Chipre
$ws.Range("B153").Value = 1490
$ws.Range("C153").Value = 2
$ws.Range("E153").Value = 330
$ws.Range("G153").Value = 1
$ws.Range("H153").Value = 21

# Row 154: Togo
$ws.Range("A154").Value = 'Togo'
$ws.Range("B154").Value = 1416
$ws.Range("C154").Value = 16
$ws.Range("D154").Value = 1035
$ws.Range("E154").Value = 353
$ws.Range("H154").Value = 28

# Row 155: Letonia
$ws.Range("A155").Value = 'Letonia'
$ws.Range("B155").Value = 1404
$ws.Range("C155").Value = 8
$ws.Range("D155").Value = 1173
$ws.Range("E155").Value = 197
$ws.Range("H155").Value = 34

# Row 156: Guyana
$ws.Range("A156").Value = 'Guyana'
$ws.Range("B156").Value = 1373
$ws.Range("C156").Value = 67
$ws.Range("D156").Value = 742
$ws.Range("E156").Value = 590
$ws.Range("G156").Value = 2
$ws.Range("H156").Value = 41

# Row 157: Burkina Faso
$ws.Range("A157").Value = 'Burkina Faso'
$ws.Range("B157").Value = 1370
$ws.Range("C157").Value = 2
$ws.Range("D157").Value = 1075
$ws.Range("E157").Value = 240
$ws.Range("H157").Value = 55

# Row 158: Liberia
$ws.Range("B158").Value = 1305
$ws.Range("C158").Value = 1
$ws.Range("D158").Value = 1158
$ws.Range("E158").Value = 65

# Row 160: Principado de Andorra
$ws.Range("B160").Value = 1184
$ws.Range("C160").Value = 8
$ws.Range("E160").Value = 223

# Row 164: Republica del Chad
$ws.Range("B164").Value = 1017
$ws.Range("C164").Value = 4
$ws.Range("D164").Value = 884

# Row 168: Martinica
$ws.Range("A168").Value = 'Martinica'
$ws.Range("B168").Value = 747
$ws.Range("C168").Value = 132
$ws.Range("D168").Value = 98
$ws.Range("E168").Value = 631
$ws.Range("G168").Value = 2
$ws.Range("H168").Value = 18

# Row 169: San Marino
$ws.Range("A169").Value = 'San Marino'
$ws.Range("B169").Value = 715
$ws.Range("D169").Value = 660
$ws.Range("E169").Value = 13
$ws.Range("H169").Value = 42

# Row 170: Crucero
$ws.Range("A170").Value = 'Crucero'
$ws.Range("B170").Value = 712
$ws.Range("D170").Value = 651
$ws.Range("E170").Value = 48
$ws.Range("H170").Value = 13

# Row 178: Comoras
$ws.Range("B178").Value = 427
$ws.Range("C178").Value = 4
$ws.Range("D178").Value = 410
$ws.Range("E178").Value = 10

# Row 190: Butan
$ws.Range("D190").Value = 157
$ws.Range("E190").Value = 6

# Row 192: Monaco
$ws.Range("B192").Value = 140
$ws.Range("C192").Value = 2
$ws.Range("D192").Value = 90
$ws.Range("E192").Value = 49
